# Update the "Data" worksheet: replace the second data row's values
# (Households and NPISHs final consumption expenditure series) with the
# revised figures, and move the active selection to B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$values = @(
    17285.148130000001,
    18075.937089999999,
    18651.36636,
    19262.448639999999,
    19503.88221,
    19257.311450000001,
    19333.654729999998,
    19430.545020000001,
    20343.653160000002,
    21229.28023,
    22131.44112,
    22833.060249999999,
    23396.543470000001,
    24153.47997,
    24624.04567,
    24837.821459999999,
    24552.37283,
    25105.307919999999,
    25639.547989999999,
    26308.300220000001,
    26763.924729999999,
    27372.299060000001,
    28063.0929,
    29211.033449999999,
    30431.437699999999,
    31610.511210000001,
    32074.27982,
    32586.169379999999,
    33329.568169999999,
    34267.479729999999,
    35155.495020000002,
    35826.942620000002,
    36350.502229999998,
    36062.319150000003,
    35280.877110000001,
    35656.258959999999,
    36000.369180000002,
    36225.887690000003,
    36603.181080000002,
    37364.14086,
    38340.950839999998,
    38998.402800000003,
    39774.807489999999,
    40650.48098,
    41334.696530000001,
    39898.878940000002,
    43330.52173
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2   # column B is index 2
    $ws.Cells.Item(2, $col).Value = $values[$i]
}

# Move the active selection from C10 to B13 on the Data sheet.
$ws.Range("B13").Select()
